# Macroferia Regional de Talca - Betarraga: insert a new daily price report.
# Row 125 ("today") keeps its data but its date advances by one day; every
# row below it (126..249) shifts down by one row, and a brand-new row 250
# is created holding what used to be the last row's data.
#
# Copying row 125 and inserting the copy at row 126 pushes rows 126..249
# down to 127..250 automatically (Excel's normal insert-row shift
# behaviour), which reproduces that whole cascade in one operation and
# leaves row 125's own values (except the date) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(125).Copy()
$ws.Rows(126).Insert()

# Bump row 125's date forward by one day (44586 -> 44587).
$ws.Range("D125").Value = $ws.Range("D125").Value2() + 1
